$d = $word.ActiveDocument
$find = $d.Content.Find

function ReplaceText($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $old"
    }
}

# 1. " (Cess)" -- merge run split by spellStart/spellEnd proofErr
ReplaceText " (Cess)" " (Cess)"

# 2. "Once risks are identified...each risk in order to allow..." -- merge run split by gramStart/gramEnd around "in order to"
ReplaceText "each risk in order to allow" "each risk in order to allow"

# 3. " (Rark)" -- merge run split by spellStart/spellEnd proofErr
ReplaceText " (Rark)" " (Rark)"

# 4. "Every project must maintain a risk register in order to track risks" -- merge run split by gramStart/gramEnd
ReplaceText "risk register in order to track" "risk register in order to track"

# 5. "...each risk, its likelihood and potential consequences..." -- merge run split by gramStart/gramEnd
ReplaceText "risk, its likelihood and potential" "risk, its likelihood and potential"

# 6. "...Team Developmentality can mitigate..." -- merge run split by spellStart/spellEnd
ReplaceText "Team Developmentality can mitigate" "Team Developmentality can mitigate"

# 7. "Probability - ... with 1 " -- merge run split by gramStart/gramEnd
ReplaceText "Probability - likelihood of a risk occurring is assessed using a scale of 1 to 5, with 1 " "Probability - likelihood of a risk occurring is assessed using a scale of 1 to 5, with 1 "

# 8. "Impact - ... with 1 " -- merge run split by gramStart/gramEnd
# (search starts just after the apostrophe in "risk's" so AutoCorrect's
# smart-quote substitution never touches that character)
ReplaceText "potential impact on the project is rated on a scale of 1 to 5, with 1 " "potential impact on the project is rated on a scale of 1 to 5, with 1 "

# 9. "Risk Score - ... overall risk " -- merge run split by gramStart/gramEnd
ReplaceText "Risk Score - the probability and impact scores are multiplied to determine the overall risk " "Risk Score - the probability and impact scores are multiplied to determine the overall risk "

# 10. "Status -risk's current status, whether..." -- merge run split by gramStart/gramEnd
# (search starts with the space right after "risk's" so the apostrophe
# itself is never part of the replaced span)
ReplaceText " current status, whether it is open, in progress, or closed, is also documented." " current status, whether it is open, in progress, or closed, is also documented."

Write-Host "All replacements done"
